$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column I (AIRCRAFT is H, CONTROL-DIM-1 was I) so the
# new "THRESHOLD-USED" column lands between AIRCRAFT and CONTROL-DIM-1.
$ws.Columns.Item(9).Insert()

# Insert a new data row before current row 3 ("1/2/3/v/d/d ... DT / 01.02.2014"),
# duplicating row 2 into it so it starts out identical to the row above
# ("1/2/3/v/d/d ... H / 10000"), then only the differing cells are overwritten
# below. This keeps cell formatting/type identical to the source data instead
# of re-typing every value (which would needlessly mark text-like numbers
# with a quote-prefix style).
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(2).Copy()
$ws.Rows.Item(3).PasteSpecial()
$excel.CutCopyMode = $false

# Header row
$ws.Cells.Item(1, 9).Value = "THRESHOLD-USED"

# Data rows - column I (THRESHOLD-USED) values
$ws.Cells.Item(2, 9).Value = "Y"
$ws.Cells.Item(3, 9).Value = "Y"
$ws.Cells.Item(4, 9).Value = "Y"
$ws.Cells.Item(5, 9).Value = "Y"
$ws.Cells.Item(6, 9).Value = "N"
$ws.Cells.Item(7, 9).Value = "N"
$ws.Cells.Item(8, 9).Value = "N"

# The newly inserted row 3 was a copy of row 2 ("... H / 10000"); fix up its
# CONTROL-DIM-1 / DUE-AMOUNT-1 cells to the new row's actual values ("C / 4000").
$ws.Cells.Item(3, 10).Value = "C"
$ws.Cells.Item(3, 11).Value = "'4000"
